# Fruta / hortaliza, semanal
# Insert a new weekly record at row 41 (shifting the existing rows 41-128
# down to 42-129) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows down by inserting a new row at position 41.
$ws.Rows.Item(41).Insert()

# Fill in the new row with the new weekly observation.
$ws.Range("A41").Value = 6
$ws.Range("B41").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C41").Value = "Metropolitana"
$ws.Range("D41").Value = 44804
$ws.Range("E41").Value = 13
$ws.Range("F41").Value = "Fruta"
$ws.Range("G41").Value = 100104
$ws.Range("H41").Value = "Frutos de pepita"
$ws.Range("I41").Value = 100104003
$ws.Range("J41").Value = "Membrillo"
$ws.Range("K41").Value = "Champion"
$ws.Range("L41").Value = "Especial"
$ws.Range("M41").Value = 7
$ws.Range("N41").Value = 250000
$ws.Range("O41").Value = 250000
$ws.Range("P41").Value = 250000
$ws.Range("Q41").Value = "$/bins (450 kilos)"
$ws.Range("R41").Value = "Región de O'Higgins"
$ws.Range("S41").Value = 556
$ws.Range("T41").Value = 450
